# Insert a new weekly price-report row for "Ají" (Vega Monumental Concepción)
# at row 45, pushing the existing rows 45-112 down to 46-113.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 45..112 down to 46..113, leaving a blank row 45 to fill in.
$ws.Rows(45).Insert()

# Populate the newly inserted row 45 with the new record.
$ws.Range("A45").Value = 11
$ws.Range("B45").Value = "Vega Monumental Concepción"
$ws.Range("C45").Value = "Bíobío"
$ws.Range("D45").Value = 44721
$ws.Range("E45").Value = 8
$ws.Range("F45").Value = 100112021
$ws.Range("G45").Value = "Ají"
$ws.Range("H45").Value = "Inferno"
$ws.Range("I45").Value = "Primera"
$ws.Range("J45").Value = 35
$ws.Range("K45").Value = 22000
$ws.Range("L45").Value = 23000
$ws.Range("M45").Value = 22571
$ws.Range("N45").Value = "$/caja 15 kilos"
$ws.Range("O45").Value = "Provincia de Huasco"
$ws.Range("P45").Value = 1505
$ws.Range("Q45").Value = 15
$ws.Range("R45").Value = "Hortaliza"
